$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.447150491076314
$ws.Range("D2").Value = 0.2099720944267673
$ws.Range("E2").Value = 0.2129104482259407
$ws.Range("F2").Value = 1.427864945571685
$ws.Range("G2").Value = 0.00244283080509517
$ws.Range("J2").Value = 0.2688436694919218
$ws.Range("L2").Value = 0.4333223554745302
$ws.Range("N2").Value = 1.594904105419729
$ws.Range("O2").Value = 3.515032203701594

$ws.Range("B3").Value = 1.384924703245559
$ws.Range("D3").Value = 0.2105145597323883
$ws.Range("E3").Value = 0.2124355301268874
$ws.Range("F3").Value = 1.425083920606397
$ws.Range("G3").Value = 0.002446311811231266
$ws.Range("J3").Value = 0.2668291431689696
$ws.Range("L3").Value = 0.3935656026381196
$ws.Range("N3").Value = 1.593305147813595
$ws.Range("O3").Value = 3.484609906073871

$ws.Range("B4").Value = 1.347145246703207
$ws.Range("D4").Value = 0.2108792011133112
$ws.Range("E4").Value = 0.2121762828696312
$ws.Range("F4").Value = 1.424255684989873
$ws.Range("G4").Value = 0.002448565787506926
$ws.Range("J4").Value = 0.265631363953446
$ws.Range("L4").Value = 0.3691657280409117
$ws.Range("N4").Value = 1.59295104811001
$ws.Range("O4").Value = 3.468234449352536

$ws.Range("B5").Value = 1.331858482527849
$ws.Range("D5").Value = 0.2110357602785378
$ws.Range("E5").Value = 0.2120788382748771
$ws.Range("F5").Value = 1.424139126451621
$ws.Range("G5").Value = 0.002449513718111685
$ws.Range("J5").Value = 0.2651532271011448
$ws.Range("L5").Value = 0.3592259882106816
$ws.Range("N5").Value = 1.592965146043042
$ws.Range("O5").Value = 3.462139993747854

$ws.Range("B6").Value = 1.329326725577488
$ws.Range("D6").Value = 0.2110622386742431
$ws.Range("E6").Value = 0.212063154892121
$ws.Range("F6").Value = 1.42413311246672
$ws.Range("G6").Value = 0.002449672900750111
$ws.Range("J6").Value = 0.2650744385205286
$ws.Range("L6").Value = 0.3575757271854059
$ws.Range("N6").Value = 1.592977070504446
$ws.Range("O6").Value = 3.461162945088063

$ws.Range("B7").Value = 1.346938642473532
$ws.Range("D7").Value = 0.2108812802380235
$ws.Range("E7").Value = 0.21217493540642
$ws.Range("F7").Value = 1.42425321861969
$ws.Range("G7").Value = 0.002448578452465335
$ws.Range("J7").Value = 0.2656248751079673
$ws.Range("L7").Value = 0.3690316625636854
$ws.Range("N7").Value = 1.592950596145371
$ws.Range("O7").Value = 3.468149915257442

$ws.Range("B8").Value = 1.425607018151965
$ws.Range("D8").Value = 0.2101526014614237
$ws.Range("E8").Value = 0.212740019928491
$ws.Range("F8").Value = 1.426723443138002
$ws.Range("G8").Value = 0.002444006905362219
$ws.Range("J8").Value = 0.2681410089539469
$ws.Range("L8").Value = 0.4196123399335931
$ws.Range("N8").Value = 1.594222844605881
$ws.Range("O8").Value = 3.50406398384817

$ws.Range("B9").Value = 1.583223319965157
$ws.Range("D9").Value = 0.2089729257416977
$ws.Range("E9").Value = 0.2141021385653303
$ws.Range("F9").Value = 1.438554206058299
$ws.Range("G9").Value = 0.002435963346704568
$ws.Range("J9").Value = 0.27338065050116
$ws.Range("L9").Value = 0.5188636461713827
$ws.Range("N9").Value = 1.601675226018074
$ws.Range("O9").Value = 3.592811637148941

$ws.Range("B10").Value = 1.701019202634711
$ws.Range("D10").Value = 0.2082565782257078
$ws.Range("E10").Value = 0.2152542463760305
$ws.Range("F10").Value = 1.451522335549114
$ws.Range("G10").Value = 0.002430609574753557
$ws.Range("J10").Value = 0.2774100713796273
$ws.Range("L10").Value = 0.5917968481687126
$ws.Range("N10").Value = 1.610145153880026
$ws.Range("O10").Value = 3.669250693703702

$ws.Range("B11").Value = 1.755031979366606
$ws.Range("D11").Value = 0.2079630076519692
$ws.Range("E11").Value = 0.2158105120089182
$ws.Range("F11").Value = 1.458354239523359
$ws.Range("G11").Value = 0.002428293470400411
$ws.Range("J11").Value = 0.2792808908055378
$ws.Range("L11").Value = 0.6249741017973065
$ws.Range("N11").Value = 1.61464310198653
$ws.Range("O11").Value = 3.706480025791336

$ws.Range("B12").Value = 1.77554557264898
$ws.Range("D12").Value = 0.2078564576009079
$ws.Range("E12").Value = 0.2160257171751425
$ws.Range("F12").Value = 1.461075663649083
$ws.Range("G12").Value = 0.002427433491301634
$ws.Range("J12").Value = 0.2799946396079491
$ws.Range("L12").Value = 0.6375368377299537
$ws.Range("N12").Value = 1.61643860885043
$ws.Range("O12").Value = 3.720932079586362

$ws.Range("B13").Value = 1.771124949020191
$ws.Range("D13").Value = 0.2078792000507796
$ws.Range("E13").Value = 0.215979167157915
$ws.Range("F13").Value = 1.460483578410987
$ws.Range("G13").Value = 0.002427617945117173
$ws.Range("J13").Value = 0.2798406868950352
$ws.Range("L13").Value = 0.6348312732718
$ws.Range("N13").Value = 1.61604782031381
$ws.Range("O13").Value = 3.71780380875623

$ws.Range("B14").Value = 1.756718447153048
$ws.Range("D14").Value = 0.2079541492959613
$ws.Range("E14").Value = 0.2158281260665653
$ws.Range("F14").Value = 1.458575439291991
$ws.Range("G14").Value = 0.002428222377546458
$ws.Range("J14").Value = 0.2793395056571484
$ws.Range("L14").Value = 0.6260076651246891
$ws.Range("N14").Value = 1.614788973845208
$ws.Range("O14").Value = 3.707661903256337

$ws.Range("B15").Value = 1.747901842009185
$ws.Range("D15").Value = 0.2080006586054814
$ws.Range("E15").Value = 0.2157362008565791
$ws.Range("F15").Value = 1.457424149248595
$ws.Range("G15").Value = 0.002428594831751179
$ws.Range("J15").Value = 0.2790332054237581
$ws.Range("L15").Value = 0.6206028290962138
$ws.Range("N15").Value = 1.614029890339552
$ws.Range("O15").Value = 3.701495833875981

$ws.Range("B16").Value = 1.69749782261016
$ws.Range("D16").Value = 0.2082764114855173
$ws.Range("E16").Value = 0.2152185347246629
$ws.Range("F16").Value = 1.451094640249394
$ws.Range("G16").Value = 0.002430763332330971
$ws.Range("J16").Value = 0.2772885589466512
$ws.Range("L16").Value = 0.5896285674024
$ws.Range("N16").Value = 1.609864136390456
$ws.Range("O16").Value = 3.666867177396512

$ws.Range("B17").Value = 1.666685020688476
$ws.Range("D17").Value = 0.2084538311520667
$ws.Range("E17").Value = 0.2149091547774837
$ws.Range("F17").Value = 1.447450717396009
$ws.Range("G17").Value = 0.00243212414684868
$ws.Range("J17").Value = 0.2762278692957807
$ws.Range("L17").Value = 0.5706262469081196
$ws.Range("N17").Value = 1.607473399377412
$ws.Range("O17").Value = 3.646253456530701

$ws.Range("B18").Value = 1.649002562606086
$ws.Range("D18").Value = 0.2085589199718179
$ws.Range("E18").Value = 0.2147342391632243
$ws.Range("F18").Value = 1.445442600456062
$ws.Range("G18").Value = 0.002432918090508618
$ws.Range("J18").Value = 0.2756213556614142
$ws.Range("L18").Value = 0.5596966168747599
$ws.Range("N18").Value = 1.606159044918996
$ws.Range("O18").Value = 3.634628187746671

$ws.Range("B19").Value = 1.643022537438526
$ws.Range("D19").Value = 0.2085950244495791
$ws.Range("E19").Value = 0.2146755384273042
$ws.Range("F19").Value = 1.444777754695437
$ws.Range("G19").Value = 0.002433188838984858
$ws.Range("J19").Value = 0.27541661710206
$ws.Range("L19").Value = 0.5559960539470694
$ws.Range("N19").Value = 1.605724473920745
$ws.Range("O19").Value = 3.630731754879264

$ws.Range("B20").Value = 1.669960939346026
$ws.Range("D20").Value = 0.2084346299083037
$ws.Range("E20").Value = 0.2149417755263094
$ws.Range("F20").Value = 1.447829533640203
$ws.Range("G20").Value = 0.002431978123173312
$ws.Range("J20").Value = 0.2763404132758396
$ws.Range("L20").Value = 0.5726490812877785
$ws.Range("N20").Value = 1.607721615364881
$ws.Range("O20").Value = 3.648423888167542

$ws.Range("B21").Value = 1.760948364095157
$ws.Range("D21").Value = 0.207932009739153
$ws.Range("E21").Value = 0.2158723672791787
$ws.Range("F21").Value = 1.459132258575892
$ws.Range("G21").Value = 0.002428044378274376
$ws.Range("J21").Value = 0.2794865716372215
$ws.Range("L21").Value = 0.6285993992570411
$ws.Range("N21").Value = 1.615156228445329
$ws.Range("O21").Value = 3.710631209614348

$ws.Range("B22").Value = 1.82076379429742
$ws.Range("D22").Value = 0.2076304281645669
$ws.Range("E22").Value = 0.2165071082586074
$ws.Range("F22").Value = 1.467302324441576
$ws.Range("G22").Value = 0.002425572962814625
$ws.Range("J22").Value = 0.2815736602863836
$ws.Range("L22").Value = 0.6651613514380585
$ws.Range("N22").Value = 1.620552543829376
$ws.Range("O22").Value = 3.753351719935154

$ws.Range("B23").Value = 1.788807571105963
$ws.Range("D23").Value = 0.2077889341808259
$ws.Range("E23").Value = 0.2161659278800805
$ws.Range("F23").Value = 1.462870082273497
$ws.Range("G23").Value = 0.002426882925547627
$ws.Range("J23").Value = 0.2804569572445601
$ws.Range("L23").Value = 0.6456482235075498
$ws.Range("N23").Value = 1.617623425702178
$ws.Range("O23").Value = 3.730361811048283

$ws.Range("B24").Value = 1.668479796017266
$ws.Range("D24").Value = 0.2084433011712719
$ws.Range("E24").Value = 0.21492701848927
$ws.Range("F24").Value = 1.447658000380187
$ws.Range("G24").Value = 0.002432044104466039
$ws.Range("J24").Value = 0.2762895218823189
$ws.Range("L24").Value = 0.5717345732097669
$ws.Range("N24").Value = 1.607609209630368
$ws.Range("O24").Value = 3.647441932613901

$ws.Range("B25").Value = 1.540230272110762
$ws.Range("D25").Value = 0.2092655397981389
$ws.Range("E25").Value = 0.2137068185770978
$ws.Range("F25").Value = 1.434603965530513
$ws.Range("G25").Value = 0.002438041317575673
$ws.Range("J25").Value = 0.2719311438723722
$ws.Range("L25").Value = 0.492009413382533
$ws.Range("N25").Value = 1.599131309853178
$ws.Range("O25").Value = 3.566834749937783
